$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "D" column labels (rows 12-27), in top-to-bottom order ---
$ws.Range("D12").Value = "DNA Sensitivity"
$ws.Range("D13").Value = "DNA specificity"
$ws.Range("D14").Value = "DNA Predictive ACC"
$ws.Range("D15").Value = "DNA MCC"
$ws.Range("D16").Value = "RNA Sensitivity"
$ws.Range("D17").Value = "RNA specificity"
$ws.Range("D18").Value = "RNA Predictive ACC"
$ws.Range("D19").Value = "RNA MCC"
$ws.Range("D20").Value = "DRNA Sensitivity"
$ws.Range("D21").Value = "DRNA specificity"
$ws.Range("D22").Value = "DRNA Predictive ACC"
$ws.Range("D23").Value = "DRNA MCC"
$ws.Range("D24").Value = "nonDNA Sensitivity"
$ws.Range("D25").Value = "nonDNA specificity"
$ws.Range("D26").Value = "nonDNA Predictive ACC"
$ws.Range("D27").Value = "nonDNA MCC"

# These two reuse already-existing shared strings ("averageMCC" / "accuracy")
$ws.Range("D28").Value = "averageMCC"
$ws.Range("D29").Value = "accuracy"

# --- Formulas that mirror the values from the table above into the new column E ---
$ws.Range("E12").Formula = "=E3"
$ws.Range("E13").Formula = "=E4"
$ws.Range("E14").Formula = "=E5"
$ws.Range("E15").Formula = "=E6"
$ws.Range("E16").Formula = "=F3"
$ws.Range("E17").Formula = "=F4"
$ws.Range("E18").Formula = "=F5"
$ws.Range("E19").Formula = "=F6"
$ws.Range("E20").Formula = "=G3"
$ws.Range("E21").Formula = "=G4"
$ws.Range("E22").Formula = "=G5"
$ws.Range("E23").Formula = "=G6"
$ws.Range("E24").Formula = "=H3"
$ws.Range("E25").Formula = "=H4"
$ws.Range("E26").Formula = "=H5"
$ws.Range("E27").Formula = "=H6"
$ws.Range("E28").Formula = "=E7"
$ws.Range("E29").Formula = "=E8"

# --- Style the D column labels with a small Arial font ---
# (build the font once on D12, then copy/paste-special just the formatting
#  onto the rest of the range so every cell ends up sharing the same style
#  without disturbing the text already entered in each cell)
$ws.Range("D12").Font.Name = "Arial"
$ws.Range("D12").Font.Size = 10
$ws.Range("D12").Font.Color = 0
$ws.Range("D12").Copy()
$ws.Range("D13:D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Add the instructional note above the new table, bold/underlined ---
$ws.Range("E11").Value = "Copy the VALUES over to the google doc. Title the design with something descriptive (ie PAAC data using Random Forest)"
$ws.Range("E11").Font.Bold = $true
$ws.Range("E11").Font.Underline = $true
$ws.Range("E11").Font.Name = "Calibri (Body)"

# --- Update the active selection to match the saved view ---
$ws.Range("G19").Select()
